# Benchmark update: 2025-11-22 06:34:55 UTC
# Clears out a set of now-unavailable bank fee figures (the rows that had
# identical boilerplate "30,46 TL - 60,94 TL - 609,43 TL"-style figures
# copy/pasted across many banks) and blanks a handful of now-unconfirmed
# "Asgari Tutar / Azami Tutar" combo figures while keeping the label text
# shell around the two remaining blank placeholders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose benchmark figure is removed entirely (cleared to empty).
$cellsToClear = @(
    "D2","I2","J2",
    "C3","D3","F3","H3","I3","J3","K3",
    "C4","D4","F4","H4","I4","J4","K4",
    "C5","D5","F5","H5","I5","J5","K5",
    "C6","D6","H6","I6","J6","K6",
    "D7","J7",
    "C8","D8","F8","H8","I8","J8","K8",
    "C9","D9","F9","H9","I9","J9","K9",
    "C10","D10","F10","H10","I10","J10","K10",
    "C11","D11","H11","I11","J11","K11",
    "C12","D12","K12",
    "C13","D13","F13","H13","I13","J13","K13",
    "C14","D14","F14","H14","J14","K14",
    "K20","K23","K24","K25"
)

foreach ($addr in $cellsToClear) {
    $ws.Range($addr).Value = ""
}

# Cells that keep the "Asgari Tutar / Azami Tutar" label template but with
# the percentage and amounts scrubbed out.
$ws.Range("K15").Value = " Asgari Tutar:  Azami Tutar: "
$ws.Range("K17").Value = " Asgari Tutar:  Azami Tutar: "
$ws.Range("K21").Value = " Asgari Tutar:  Azami Tutar: "
$ws.Range("K22").Value = " Asgari Tutar:  Azami Tutar: "
